$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a header row (row 1) followed by 20 data rows
# (rows 2-21). Two new rows of data were added right after the header,
# pushing the existing data down by two rows, and eight more rows of
# data were appended after what is now the last existing row.

# Insert two blank rows right after the header row; this shifts the
# existing data (originally rows 2-21) down to rows 4-23.
$ws.Rows.Item(2).Resize(2).Insert()

# Excel's row insert copies formatting from the row above (the bold
# header style) onto the freshly inserted rows; clear that so the new
# rows look like ordinary, unstyled data rows.
$ws.Range("A2:C3").ClearFormats()

# Populate the two newly inserted rows.
$ws.Cells.Item(2, 1).Value = 0.0100699262883591
$ws.Cells.Item(2, 2).Value = -0.003861541194575109
$ws.Cells.Item(2, 3).Value = -0.02743906991518268

$ws.Cells.Item(3, 1).Value = -0.04216528505238931
$ws.Cells.Item(3, 2).Value = -0.05587235412427344
$ws.Cells.Item(3, 3).Value = -0.005946585338334982

# Append eight new rows of data after the existing data (which, after
# the insert above, now ends at row 23).
$ws.Cells.Item(24, 1).Value = 1.138076220239913
$ws.Cells.Item(24, 2).Value = 4.945667840996566
$ws.Cells.Item(24, 3).Value = -0.5911586260309147

$ws.Cells.Item(25, 1).Value = 0.5868015289306701
$ws.Cells.Item(25, 2).Value = 4.023616756711703
$ws.Cells.Item(25, 3).Value = 0.7568838426044971

$ws.Cells.Item(26, 1).Value = -0.3439888250538894
$ws.Cells.Item(26, 2).Value = 1.417871174155451
$ws.Cells.Item(26, 3).Value = 1.152574896812441

$ws.Cells.Item(27, 1).Value = -0.006142936684953748
$ws.Cells.Item(27, 2).Value = 0.247877272416142
$ws.Cells.Item(27, 3).Value = -0.5417658090591317

$ws.Cells.Item(28, 1).Value = -0.1216962014001841
$ws.Cells.Item(28, 2).Value = -0.5174121899264155
$ws.Cells.Item(28, 3).Value = -0.2046180449578262

$ws.Cells.Item(29, 1).Value = -0.007408298704090516
$ws.Cells.Item(29, 2).Value = -0.9081197368855365
$ws.Cells.Item(29, 3).Value = -0.06973525623277696

$ws.Cells.Item(30, 1).Value = 0.06267290592801897
$ws.Cells.Item(30, 2).Value = -0.8978846316434917
$ws.Cells.Item(30, 3).Value = -0.05514305708359769

$ws.Cells.Item(31, 1).Value = -0.03725966301803718
$ws.Cells.Item(31, 2).Value = -0.5934867311497137
$ws.Cells.Item(31, 3).Value = 0.0119181060973484
